$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap the match data (columns F:V) between row 128 and row 129 ---
# Columns A:E (Indice, pais, torneio, temporada, data_partida) stay put on each row;
# only the match-specific columns F..V trade places between the two rows.
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$row128vals = @{}
$row129vals = @{}
foreach ($col in $cols) {
    $row128vals[$col] = $ws.Range("$col" + "128").Value2
    $row129vals[$col] = $ws.Range("$col" + "129").Value2
}
foreach ($col in $cols) {
    $ws.Range("$col" + "128").Value2 = $row129vals[$col]
    $ws.Range("$col" + "129").Value2 = $row128vals[$col]
}

# --- Step 2: append a new row 130 with the new match ---
# Copy formatting (styles) from row 129 down to the new row 130 first.
$ws.Range("A129:V129").Copy() | Out-Null
$ws.Range("A130:V130").PasteSpecial(-4122) | Out-Null

$ws.Range("A130").Value2 = 129
$ws.Range("B130").Value2 = "paraguay"
$ws.Range("C130").Value2 = "primera-division"
# D130 ("2023") must stay text, not auto-coerced to a number - force text format,
# assign, then restore the default "Normal" style so no stray numFmt is left behind.
$ws.Range("D130").NumberFormat = "@"
$ws.Range("D130").Value2 = "2023"
$ws.Range("D130").Style = "Normal"
$ws.Range("E130").Value2 = 45261.02083333334
$ws.Range("F130").Value2 = "Resistencia"
$ws.Range("G130").Value2 = 0
$ws.Range("H130").Value2 = "Ameliano"
$ws.Range("I130").Value2 = 3
$ws.Range("J130").Value2 = 3.79
$ws.Range("K130").Value2 = "24/11/2023 23:42"
$ws.Range("L130").Value2 = 7.17
$ws.Range("M130").Value2 = "01/12/2023 00:29"
$ws.Range("N130").Value2 = 3.65
$ws.Range("O130").Value2 = "24/11/2023 23:42"
$ws.Range("P130").Value2 = 4.94
$ws.Range("Q130").Value2 = "01/12/2023 00:29"
$ws.Range("R130").Value2 = 1.97
$ws.Range("S130").Value2 = "24/11/2023 23:42"
$ws.Range("T130").Value2 = 1.44
$ws.Range("U130").Value2 = "01/12/2023 00:27"
$ws.Range("V130").Value2 = "https://www.betexplorer.com/football/paraguay/primera-division/resistencia-sportivo-ameliano/ER8Ao38t/"
